$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert "franzosa_ControlvsCD_ConvCD" row ---
# Currently row 8 is "franzosa_ControlvsCD_Fp"; push it (and everything below) down
# by inserting a new blank row at row 8, then fill it in.
$ws.Rows.Item(8).Insert()

$ws.Cells.Item(8, 1).Value = "franzosa_ControlvsCD_ConvCD"
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1
$ws.Cells.Item(8, 8).Value = 1

# --- Insert "franzosa_ControlvsUC_ConvUC" row ---
# After the previous insertion, "franzosa_ControlvsUC_Age" is now row 13 and
# "franzosa_ControlvsUC_Fp" is now row 14. Insert a new blank row at row 14
# (before the shifted "franzosa_ControlvsUC_Fp") and fill it in.
$ws.Rows.Item(14).Insert()

$ws.Cells.Item(14, 1).Value = "franzosa_ControlvsUC_ConvUC"
$ws.Cells.Item(14, 2).Value = 0
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 8).Value = 0
